# MFP_Historical_Reservoir_Data.xlsx - adjust the "Historical Snowpack" (D) column
# so the spreadsheet's reference water-year doesn't need to be edited every year.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# D15 (2020-10-14): 0 -> 1
$ws.Range("D15").Value = 1

# D319:D358 (2021-06-26 through 2021-08-04): set to a flat 0.01 baseline,
# replacing the previous all-zero (with one stray 0.158) placeholder data.
$ws.Range("D319:D358").Value = 0.01

# Leave the selection on the range that was just edited, scrolled into view.
$ws.Range("D319:D358").Select()
$excel.ActiveWindow.ScrollRow = 322
